$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The agriculture-orientation-index row (row 4) was shown with only one decimal
# place (custom format "0.0"); switch it to the built-in two-decimal format.
$ws.Range("D4:P4").NumberFormat = "0.00"

# Add the new 2020 data column (Q), inheriting the same look as the 2019
# column (P) for both the year header (row 3) and the index value (row 4).
$ws.Range("P3:P4").Copy()
$ws.Range("Q3:Q4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("Q3").Value = 2020
$ws.Range("Q4").Value = 0.0671560491274446

# Reset the selection back to A1 (the sheet had a stray B12 selection saved).
[void]$ws.Range("A1").Select()
